$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (chars). The host quantizes ColumnWidth writes to a pixel
# grid (1/6-character steps) and then adds a constant 5/6-character offset
# again on export, so the input is pre-compensated (target - 5/6) to land
# the saved width on the closest achievable grid point to the target.
$ws.Columns.Item(2).ColumnWidth = 50.833333333333336
$ws.Columns.Item(3).ColumnWidth = 47.833333333333336
$ws.Columns.Item(4).ColumnWidth = 54.833333333333336
$ws.Columns.Item(6).ColumnWidth = 55.833333333333336

# Row 2
$ws.Range("C2").Value = "{0: sala nr 8 | Piotr Wójcik | Biologia}"
$ws.Range("E2").Value = "{}"

# Row 3
$ws.Range("C3").Value = "{0: sala nr 5 | Katarzyna Mazur | Fizyka}"
$ws.Range("E3").Value = "{}"

# Row 4
$ws.Range("C4").Value = "{0: sala nr 3 | Lena Kowalska | Język angielski}"
$ws.Range("E4").Value = "{}"

# Row 5
$ws.Range("C5").Value = "{0: sala nr 8 | Paweł Lewandowski | Matematyka}"
$ws.Range("E5").Value = "{}"

# Row 6
$ws.Range("C6").Value = "{0: sala nr 2 | Jan Nowak | Język polski}"
$ws.Range("D6").Value = "{0: sala nr 6 | Paweł Lewandowski | Matematyka}"
$ws.Range("E6").Value = "{}"

# Row 7
$ws.Range("B7").Value = "{0: sala nr 3 | Lena Kowalska | Język angielski}"
$ws.Range("C7").Value = "{}"
$ws.Range("D7").Value = "{0: sala nr 9 | Mateusz Kowalski | Język niemiecki}"
$ws.Range("E7").Value = "{}"
$ws.Range("F7").Value = "{0: sala nr 11 | Zofia Wiśniewska | Wychowanie fizyczne}"

# Row 8
$ws.Range("B8").Value = "{0: sala nr 2 | Jan Nowak | Język polski}"
$ws.Range("C8").Value = "{}"
$ws.Range("D8").Value = "{0: sala nr 9 | Natalia Szymańska | Geografia}"
$ws.Range("E8").Value = "{}"
$ws.Range("F8").Value = "{0: sala nr 7 | Natalia Szymańska | Geografia}"

# Row 9
$ws.Range("B9").Value = "{0: sala nr 4 | Dominik Kaczor | Informatyka}"
$ws.Range("C9").Value = "{}"
$ws.Range("D9").Value = "{0: sala nr 5 | Karolina Kamińska | Chemia}"
$ws.Range("E9").Value = "{}"
$ws.Range("F9").Value = "{0: sala nr 4 | Paweł Lewandowski | Matematyka}"

# Row 10
$ws.Range("B10").Value = "{0: sala nr 1 | Mateusz Kowalski | Język niemiecki}"
$ws.Range("C10").Value = "{}"
$ws.Range("D10").Value = "{0: sala nr 10 | Piotr Wójcik | Biologia}"
$ws.Range("E10").Value = "{0: sala nr 5 | Lena Kowalska | Język angielski}"
$ws.Range("F10").Value = "{0: sala nr 7 | Dominik Kaczor | Informatyka}"

# Row 11
$ws.Range("B11").Value = "{0: sala nr 1 | Dominik Kaczor | Informatyka}"
$ws.Range("C11").Value = "{}"
$ws.Range("D11").Value = "{0: sala nr 1 | Katarzyna Mazur | Fizyka}"
$ws.Range("E11").Value = "{0: sala nr 8 | Dominik Kaczor | Informatyka}"
$ws.Range("F11").Value = "{}"

# Row 12
$ws.Range("B12").Value = "{0: sala nr 6 | Jan Nowak | Język polski}"
$ws.Range("C12").Value = "{}"
$ws.Range("D12").Value = "{0: sala nr 9 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("E12").Value = "{0: sala nr 5 | Paweł Lewandowski | Matematyka}"
$ws.Range("F12").Value = "{0: sala nr 5 | Paweł Lewandowski | Matematyka}"

# Row 13
$ws.Range("B13").Value = "{0: sala nr 6 | Katarzyna Mazur | Fizyka}"
$ws.Range("C13").Value = "{}"
$ws.Range("D13").Value = "{0: sala nr 6 | Dominik Kaczor | Informatyka}"
$ws.Range("E13").Value = "{0: sala nr 6 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("F13").Value = "{0: sala nr 6 | Karolina Kamińska | Chemia}"
